$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1654.25
$ws.Range("I137").Value = 1330.4166
$ws.Range("J137").Value = 2625.75
$ws.Range("K137").Value = 3991.2498
$ws.Range("L137").Value = 7877.25
$ws.Range("M137").Value = -1441.2498
$ws.Range("N137").Value = -12977.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 846.02856
$ws.Range("I45").Value = 800.1111
$ws.Range("J45").Value = 894.64703
$ws.Range("K45").Value = 800.1111
$ws.Range("L45").Value = 894.64703
$ws.Range("M45").Value = -423.1111
$ws.Range("N45").Value = -1648.64703

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1708
$ws.Range("I61").Value = 804.4286
$ws.Range("J61").Value = 3289.25
$ws.Range("K61").Value = 804.4286
$ws.Range("L61").Value = 3289.25
$ws.Range("M61").Value = -592.4286
$ws.Range("N61").Value = -3713.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1044.2858
$ws.Range("I74").Value = 1093.8889
$ws.Range("J74").Value = 746.6667
$ws.Range("K74").Value = 1093.8889
$ws.Range("L74").Value = 746.6667
$ws.Range("M74").Value = -219.8888999999999
$ws.Range("N74").Value = -2494.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1044.2858
$ws.Range("I77").Value = 1093.8889
$ws.Range("J77").Value = 746.6667
$ws.Range("K77").Value = 5469.4445
$ws.Range("L77").Value = 3733.3335
$ws.Range("M77").Value = -1101.4445
$ws.Range("N77").Value = -12469.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 317.30435
$ws.Range("I97").Value = 289.15
$ws.Range("J97").Value = 505
$ws.Range("K97").Value = 289.15
$ws.Range("L97").Value = 505
$ws.Range("M97").Value = 206.85
$ws.Range("N97").Value = -1497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1135.5
$ws.Range("I110").Value = 801.17645
$ws.Range("J110").Value = 1652.1818
$ws.Range("K110").Value = 801.17645
$ws.Range("L110").Value = 1652.1818
$ws.Range("M110").Value = 1243.82355
$ws.Range("N110").Value = -5742.1818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1076.5
$ws.Range("I122").Value = 1022.2
$ws.Range("J122").Value = 1212.25
$ws.Range("K122").Value = 3066.6
$ws.Range("L122").Value = 3636.75
$ws.Range("M122").Value = -616.6000000000004
$ws.Range("N122").Value = -8536.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5055.2
$ws.Range("I132").Value = 5100.724
$ws.Range("K132").Value = 15302.172
$ws.Range("M132").Value = -12772.172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1708
$ws.Range("I136").Value = 804.4286
$ws.Range("J136").Value = 3289.25
$ws.Range("K136").Value = 2413.2858
$ws.Range("L136").Value = 9867.75
$ws.Range("M136").Value = 136.7142000000003
$ws.Range("N136").Value = -14967.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 41666.668
$ws.Range("J92").Value = 41666.668
$ws.Range("L92").Value = 41666.668
$ws.Range("N92").Value = -46658.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22797.854
$ws.Range("I134").Value = 27389.666
$ws.Range("J134").Value = 2900
$ws.Range("K134").Value = 82168.99800000001
$ws.Range("L134").Value = 8700
$ws.Range("M134").Value = -79633.99800000001
$ws.Range("N134").Value = -13770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 2918.182
$ws.Range("J31").Value = 20002984
$ws.Range("K31").Value = 2918.182
$ws.Range("L31").Value = 20002984
$ws.Range("M31").Value = -2623.182
$ws.Range("N31").Value = -20003574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I34").Value = 2918.182
$ws.Range("J34").Value = 20002984
$ws.Range("K34").Value = 2918.182
$ws.Range("L34").Value = 20002984
$ws.Range("M34").Value = -2716.182
$ws.Range("N34").Value = -20003388

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 986.7143
$ws.Range("I58").Value = 1001.75
$ws.Range("J58").Value = 938.6
$ws.Range("K58").Value = 1001.75
$ws.Range("L58").Value = 938.6
$ws.Range("M58").Value = -798.75
$ws.Range("N58").Value = -1344.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4601.25
$ws.Range("I132").Value = 3228
$ws.Range("J132").Value = 5974.5
$ws.Range("K132").Value = 9684
$ws.Range("L132").Value = 17923.5
$ws.Range("M132").Value = -7154
$ws.Range("N132").Value = -22983.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1383.7693
$ws.Range("I134").Value = 1332.4166
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 3997.2498
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -1462.2498
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 986.7143
$ws.Range("I136").Value = 1001.75
$ws.Range("J136").Value = 938.6
$ws.Range("K136").Value = 3005.25
$ws.Range("L136").Value = 2815.8
$ws.Range("M136").Value = -455.25
$ws.Range("N136").Value = -7915.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 957487.7
$ws.Range("J131").Value = 1323581.2
$ws.Range("L131").Value = 3970743.6
$ws.Range("N131").Value = -3980823.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2169.0852
$ws.Range("I122").Value = 1917.6957
$ws.Range("J122").Value = 2410
$ws.Range("K122").Value = 5753.0871
$ws.Range("L122").Value = 7230
$ws.Range("M122").Value = -3303.0871
$ws.Range("N122").Value = -12130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 94084.37
$ws.Range("I132").Value = 108150.42
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 324451.26
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -321921.26
$ws.Range("N132").Value = -20058.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 84938.25
$ws.Range("I40").Value = 126894.875
$ws.Range("K40").Value = 126894.875
$ws.Range("M40").Value = -126758.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4009.0588
$ws.Range("I132").Value = 3743.7334
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 11231.2002
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -8701.200199999999
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10450.333
$ws.Range("I136").Value = 14175.5
$ws.Range("K136").Value = 42526.5
$ws.Range("M136").Value = -39976.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2262.6365
$ws.Range("I122").Value = 2840.8
$ws.Range("J122").Value = 1780.8334
$ws.Range("K122").Value = 8522.400000000001
$ws.Range("L122").Value = 5342.5002
$ws.Range("M122").Value = -6072.400000000001
$ws.Range("N122").Value = -10242.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1885.8108
$ws.Range("I132").Value = 1466.8387
$ws.Range("J132").Value = 4050.5
$ws.Range("K132").Value = 4400.5161
$ws.Range("L132").Value = 12151.5
$ws.Range("M132").Value = -1870.5161
$ws.Range("N132").Value = -17211.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6923.4287
$ws.Range("I136").Value = 6923.4287
$ws.Range("K136").Value = 20770.2861
$ws.Range("M136").Value = -18220.2861
